# Modified logic for order handling vendornumber checking into name checking
#
# Adds new rows to the "Constants" sheet (order / claim business-exception
# constants) and to the "Assets" sheet (new Orchestrator assets for vendor
# numbers / names and claim cost limits), matching the author's commit.

$wb = $excel.ActiveWorkbook
$wsConst  = $wb.Worksheets.Item("Constants")
$wsAssets = $wb.Worksheets.Item("Assets")

# --- Constants sheet: new business-exception rows (26-34) -------------------
# Written in an order that mirrors the shared-string table produced by the
# original authoring session.
$wsConst.Range("B26").Value = "订单-01-订单号错误"
$wsConst.Range("A26").Value = "BE-Order-1"
$wsConst.Range("A27").Value = "BE-Order-2"
$wsConst.Range("A28").Value = "BE-Order-3"
$wsConst.Range("A29").Value = "BE-Order-4"
$wsConst.Range("B27").Value = "订单-02-非初次收货"
$wsConst.Range("B28").Value = "订单-03-商品NOF"
$wsConst.Range("B29").Value = "订单-04-收货数量大于订单"
$wsConst.Range("B31").Value = "索赔-01-商品NOF"
$wsConst.Range("A34").Value = "BE-Claim-4"
$wsConst.Range("A31").Value = "BE-Claim-1"
$wsConst.Range("A32").Value = "BE-Claim-2"
$wsConst.Range("A33").Value = "BE-Claim-3"
$wsConst.Range("B32").Value = "索赔-02-供应商号错误"
$wsConst.Range("B34").Value = "索赔-04-索赔金额错误"

# --- Assets sheet: new Orchestrator assets (22-25) --------------------------
$wsAssets.Range("B22").Value = "1128_VendorNumbers"
$wsAssets.Range("A22").Value = "VendorNumbers"

$wsConst.Range("B33").Value = "索赔-03-成本超100%/SMART当前实际索赔成本/<SmartCost>"

$wsAssets.Range("B23").Value = "1128_ClaimLowerLimit"
$wsAssets.Range("B24").Value = "1128_ClaimUpperLimit"
$wsAssets.Range("A23").Value = "ClaimLowerLimit"
$wsAssets.Range("A24").Value = "ClaimUpperLimit"

$wsConst.Range("B30").Value = "订单-05-供应商号错误"
$wsConst.Range("A30").Value = "BE-Order-5"

$wsAssets.Range("A25").Value = "VendorNames"
$wsAssets.Range("B25").Value = "1128_VendorNames"

# --- View state: selections / scroll position -------------------------------
# Assets stays the active ("tabSelected") sheet; Constants scrolls down and
# both sheets land on a new selection near the bottom of their new data.
$wsConst.Activate()
$winConst = $excel.ActiveWindow
$winConst.ScrollRow = 10
$wsConst.Range("A35").Select()

$wsAssets.Activate()
$wsAssets.Range("B26").Select()
